# Insert two new rows at 545-546 (shifting existing rows 545+ down by two),
# then populate the two new rows with the new weekly price records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("545:546").Insert()

# New row 545: Pimiento "Zafiro rojo"
$ws.Range("A545").Value = 5
$ws.Range("B545").Value = "Macroferia Regional de Talca"
$ws.Range("C545").Value = "Maule"
$ws.Range("D545").Value = 44769
$ws.Range("E545").Value = 7
$ws.Range("F545").Value = 100112002
$ws.Range("G545").Value = "Pimiento"
$ws.Range("H545").Value = "Zafiro rojo"
$ws.Range("I545").Value = "Primera"
$ws.Range("J545").Value = 300
$ws.Range("K545").Value = 25000
$ws.Range("L545").Value = 25000
$ws.Range("M545").Value = 25000
$ws.Range("N545").Value = "`$/caja 15 kilos"
$ws.Range("O545").Value = "Región de Arica y Parinacota"
$ws.Range("P545").Value = 1667
$ws.Range("Q545").Value = 15
$ws.Range("R545").Value = "Hortaliza"

# New row 546: Pimiento "Zafiro verde"
$ws.Range("A546").Value = 5
$ws.Range("B546").Value = "Macroferia Regional de Talca"
$ws.Range("C546").Value = "Maule"
$ws.Range("D546").Value = 44769
$ws.Range("E546").Value = 7
$ws.Range("F546").Value = 100112002
$ws.Range("G546").Value = "Pimiento"
$ws.Range("H546").Value = "Zafiro verde"
$ws.Range("I546").Value = "Primera"
$ws.Range("J546").Value = 300
$ws.Range("K546").Value = 18000
$ws.Range("L546").Value = 18000
$ws.Range("M546").Value = 18000
$ws.Range("N546").Value = "`$/caja 15 kilos"
$ws.Range("O546").Value = "Región de Arica y Parinacota"
$ws.Range("P546").Value = 1200
$ws.Range("Q546").Value = 15
$ws.Range("R546").Value = "Hortaliza"
